$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: force text storage for numeric-looking Price (D) values
# by temporarily switching the cell to text format, then restoring the default
# (unstyled) appearance so the XML matches the original's lack of an 's' attribute.

# --- Update Price (D) and Volume(1h) (E) columns for rows with changed values ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.521.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.728.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.26%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4790'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2669'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06220'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.731.19'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07167'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.67'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6140'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.526'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.523.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006965'
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.952.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.525'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.900'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.282'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.34'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.789'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.403'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.978'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07964'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.706'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04584'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.23%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.617'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9938'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.88%  '
$ws.Range("E37").Value = '  +1.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.407'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.57%  '
$ws.Range("E41").Value = '  -6.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.006'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01500'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.577'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3876'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.994'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.07%  '
$ws.Range("E47").Value = '  +1.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05346'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.97'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.826'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.260'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.87%  '

# --- Row 38/39: RenderToken and TrustWalletToken swapped positions ---
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.093'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.08%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9155'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.12%  '
